$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Handle the PEPE/OKB row swap (rows 35 and 36) first, including B/C/D/E columns
$ws.Range("B35").Value2 = 'OKB'
$ws.Range("C35").Value2 = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("B36").Value2 = 'PEPE'
$ws.Range("C36").Value2 = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'

# Update numeric/percentage text cells in column D (Price) and E (Volume(1h)).
# NumberFormat is forced to Text ("@") before assignment so Excel does not
# reinterpret/reformat these values (e.g. "1.20" -> 1.2, or dates, etc.).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value2 = '70.961.20'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value2 = '  -2.36%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value2 = '3.842.13'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value2 = '  -3.24%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value2 = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = '593.22'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value2 = '  +1.61%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = '165.99'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value2 = '  +4.73%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value2 = '0.671'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value2 = '  -1.10%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value2 = '  +0.31%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = '0.750'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value2 = '  +0.38%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = '0.175'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value2 = '  +4.92%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = '53.07'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value2 = '  -1.03%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = '0.0000319'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value2 = '  +0.75%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value2 = '11.13'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value2 = '  +2.85%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value2 = '4.468.05'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value2 = '  -2.97%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = '3.863.26'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value2 = '  -3.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value2 = '20.65'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value2 = '  +1.31%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = '13.78'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value2 = '  -1.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value2 = '1.20'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value2 = '  -5.50%  '
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value2 = '  -1.98%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = '70.943.29'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value2 = '  -2.03%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = '432.63'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value2 = '  +0.29%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = '4.71'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value2 = '  +0.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = '93.89'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value2 = '  -1.93%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value2 = '3.25'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value2 = '  -4.63%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = '13.70'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value2 = '  -3.70%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = '4.08'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value2 = '  -7.74%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = '10.82'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value2 = '  -4.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = '5.93'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value2 = '  +0.25%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = '10.11'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value2 = '  -6.35%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = '34.82'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value2 = '  -4.09%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = '7.80'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value2 = '  -0.46%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value2 = '50.36'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value2 = '  -0.82%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = '13.46'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value2 = '  -0.46%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = '0.124'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value2 = '  -5.08%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = '68.70'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value2 = '  +0.59%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = '0.0₃0978'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value2 = '  +14.05%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = '614.59'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value2 = '  -9.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value2 = '0.417'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value2 = '  -4.75%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value2 = '  -0.06%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value2 = '  +0.20%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = '3.26'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value2 = '  -1.51%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value2 = '0.141'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value2 = '  -2.38%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = '3.17'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value2 = '  +31.49%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = '0.0465'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value2 = '  -4.03%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value2 = '10.13'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value2 = '  -7.21%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value2 = '0.143'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value2 = '  -3.75%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = '2.61'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value2 = '  -2.41%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = '3.33'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value2 = '  -1.74%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value2 = '2.824.13'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value2 = '  +2.40%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = '2.72'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value2 = '  -19.52%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = '0.000270'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value2 = '  +0.09%  '
